$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 ("create customer_ID" test case) keeps its values, but the generated
# customer id in column D is rewritten to reflect the additional test steps
# that were appended since the last run.
$ws.Range("D3").Value = "sharonH-830-539-29-32-859-882"

# Row 4 is a brand-new test case ("delete_customerID"). It reuses the same
# layout/formatting as row 3, so clone the row first and then patch the
# handful of cells that actually differ.
$ws.Range("A3:K3").Copy($ws.Range("A4:K4"))
$ws.Range("A4").Value = "delete_customerID"
$ws.Range("D4").Value = "sharonH-830-539-29"

# Row 3 had explicit hyperlinks on B3 (mailto to the login email) and K3
# (mailto to the signup email); row 4 needs the same pair of hyperlinks.
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:mohitjoe91@gmail.com", "", "", "mohitjoe91@gmail.com")
$ws.Hyperlinks.Add($ws.Range("K4"), "mailto:xyz@gmail.com", "", "", "xyz@gmail.com")

# The author's selection ended up resting on the newly added row.
$ws.Range("A4").Select()
